$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was added for the "Perejil" (parsley)
# vegetable at "Feria Lagunitas de Puerto Montt". It is inserted as the
# new row 300, pushing the existing rows 300-337 down to 301-338.
$ws.Rows.Item(300).Insert()

$ws.Range("A300").Value = 4
$ws.Range("B300").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C300").Value = "Los Lagos"
$ws.Range("D300").Value = 44946
$ws.Range("E300").Value = 10
$ws.Range("F300").Value = 100112044
$ws.Range("G300").Value = "Perejil"
$ws.Range("H300").Value = "Sin especificar"
$ws.Range("I300").Value = "Primera"
$ws.Range("J300").Value = 180
$ws.Range("K300").Value = 6000
$ws.Range("L300").Value = 6000
$ws.Range("M300").Value = 6000
$ws.Range("N300").Value = "`$/docena de atados (2 kilos)"
$ws.Range("O300").Value = "Región de La Araucanía"
$ws.Range("P300").Value = 3000
$ws.Range("Q300").Value = 2
$ws.Range("R300").Value = "Hortaliza"
